$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 469 (existing rows 469-583 shift down to 471-585).
$ws.Range("A469:T470").Insert()

# The two newly inserted blank rows should carry the same "category" data
# (Mercado, Region, Producto, Categoria, Variedad, Calidad, Unidad, Origen, Kg/unidad)
# as the rows that are now immediately below them (471 and 472, formerly 469 and 470),
# since this is simply a new weekly price observation appended ahead of the existing
# Pinton / Primera Pinton pair for this market.
for ($col = 1; $col -le 20; $col++) {
    $ws.Cells.Item(469, $col).Value = $ws.Cells.Item(471, $col).Value()
    $ws.Cells.Item(470, $col).Value = $ws.Cells.Item(472, $col).Value()
}

# New row 469 ("Pintón") data point
$ws.Cells.Item(469, 4).Value = 44642
$ws.Cells.Item(469, 13).Value = 850
$ws.Cells.Item(469, 14).Value = 18000
$ws.Cells.Item(469, 15).Value = 18000
$ws.Cells.Item(469, 16).Value = 18000
$ws.Cells.Item(469, 19).Value = 900

# New row 470 ("Primera Pintón") data point
$ws.Cells.Item(470, 4).Value = 44642
$ws.Cells.Item(470, 13).Value = 500
$ws.Cells.Item(470, 14).Value = 19000
$ws.Cells.Item(470, 15).Value = 19000
$ws.Cells.Item(470, 16).Value = 19000
$ws.Cells.Item(470, 19).Value = 950
